# Generate Report for Handback
# - Update the "Status" text from "Ready for handoff" to
#   "Handed back: in sync with en-US" (Overview + per-language sheets).
# - Stamp the new handback datetimes onto the per-language sheets'
#   "Latest Handback DateTime" column (H).
# - Record the handback source (.md) / target (.xlf) files in the
#   previously-empty "Latest Handback File" (F) / "Latest Handback
#   DateTime" staging columns (F/G) for each data row, mirroring the
#   existing handoff hyperlinks in columns A/D.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Status text update
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $statusNew
$wsOverview.Range("C2").Value2 = $statusNew
$wsOverview.Range("B3").Value2 = $statusNew
$wsOverview.Range("C3").Value2 = $statusNew

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $statusNew
$wsZhCn.Range("C3").Value2 = $statusNew

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $statusNew
$wsDeDe.Range("C3").Value2 = $statusNew

# ---------------------------------------------------------------------
# 2) Latest Handback DateTime (column H) per language
# ---------------------------------------------------------------------
$wsZhCn.Range("H2").Value2 = "2016-03-23 09:58:57"
$wsZhCn.Range("H3").Value2 = "2016-03-23 09:58:57"

$wsDeDe.Range("H2").Value2 = "2016-03-23 09:59:11"
$wsDeDe.Range("H3").Value2 = "2016-03-23 09:59:11"

# ---------------------------------------------------------------------
# 3) Populate Latest Handback File (F) / Latest Handback DateTime-file
#    (G) columns with the handback source/target file links, reusing
#    the same targets as the existing handoff hyperlinks in A/D.
# ---------------------------------------------------------------------
function Get-HyperlinkUrl($ws, [string]$addr) {
    $result = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $result = $hl.Address
        }
    }
    return $result
}

foreach ($item in @(
        @{ ws = $wsZhCn; mdDisplay = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"; xlfDisplay = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.zh-cn.xlf" },
        @{ ws = $wsDeDe; mdDisplay = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.md"; xlfDisplay = "99cc93e9-5b66-41d4-bdee-c65e1d20c94d.8f480fbcd8818691a814bec4f3bb5bc290bac6b0.de-de.xlf" }
    )) {
    $ws = $item.ws
    $mdUrl = Get-HyperlinkUrl $ws '$A$2'
    $xlfUrl = Get-HyperlinkUrl $ws '$D$2'

    foreach ($row in @(2, 3)) {
        $fCell = "F" + $row
        $gCell = "G" + $row
        $ws.Hyperlinks.Add($ws.Range($fCell), $mdUrl, "", "", $item.mdDisplay)
        $ws.Hyperlinks.Add($ws.Range($gCell), $xlfUrl, "", "", $item.xlfDisplay)
    }
}

Write-Output "Handback report generated"
